$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-23 06:19:16'
$ws.Range("O2").Value = '3.4 °C'
$ws.Range("E3").Value = '2026-02-23 06:19:19'
$ws.Range("E4").Value = '2026-02-23 06:19:22'
$ws.Range("J4").Value = '1025.9 hPa'
$ws.Range("O4").Value = '5.4 °C'
$ws.Range("E5").Value = '2026-02-23 06:19:24'
$ws.Range("H5").Value = '''37%'
$ws.Range("E6").Value = '2026-02-23 06:19:27'
$ws.Range("O6").Value = '9.2 °C'
$ws.Range("E7").Value = '2026-02-23 06:19:29'
$ws.Range("J7").Value = '1025.0 hPa'
$ws.Range("N7").Value = '11.1 °C 5:58 TU'
$ws.Range("O7").Value = '11.8 °C'
$ws.Range("E8").Value = '2026-02-23 06:19:32'
$ws.Range("O8").Value = '12.9 °C'
$ws.Range("E9").Value = '2026-02-23 06:19:35'
$ws.Range("O9").Value = '6.9 °C'
$ws.Range("E10").Value = '2026-02-23 06:19:37'
$ws.Range("N10").Value = '3.2 °C 5:39 TU'
$ws.Range("O10").Value = '4.1 °C'
$ws.Range("E11").Value = '2026-02-23 06:19:40'
$ws.Range("N11").Value = '1.8 °C 5:59 TU'
$ws.Range("E12").Value = '2026-02-23 06:19:42'
$ws.Range("N12").Value = '3.1 °C 5:59 TU'
$ws.Range("O12").Value = '5.4 °C'
$ws.Range("E13").Value = '2026-02-23 06:19:45'
$ws.Range("N13").Value = '-2.8 °C 5:59 TU'
$ws.Range("E14").Value = '2026-02-23 06:19:48'
$ws.Range("H14").Value = '''84%'
$ws.Range("L14").Value = '23.4 km/h - 323º 5:57 TU'
$ws.Range("M14").Value = '9.6 °C 5:59 TU'
$ws.Range("O14").Value = '8.7 °C'
$ws.Range("E15").Value = '2026-02-23 06:19:50'
$ws.Range("O15").Value = '6.7 °C'
$ws.Range("E16").Value = '2026-02-23 06:19:52'
$ws.Range("H16").Value = '''18%'
$ws.Range("O16").Value = '2.8 °C'
$ws.Range("E17").Value = '2026-02-23 06:19:54'
$ws.Range("H17").Value = '''46%'
$ws.Range("E18").Value = '2026-02-23 06:19:57'
$ws.Range("N18").Value = '1.7 °C 5:50 TU'
$ws.Range("O18").Value = '2.9 °C'
$ws.Range("E19").Value = '2026-02-23 06:19:59'
$ws.Range("O19").Value = '9.8 °C'
$ws.Range("E20").Value = '2026-02-23 06:20:02'
$ws.Range("L20").Value = '24.1 km/h - 323º 5:40 TU'
$ws.Range("E21").Value = '2026-02-23 06:20:05'
$ws.Range("H21").Value = '''81%'
$ws.Range("N21").Value = '1.3 °C 5:40 TU'
$ws.Range("O21").Value = '3.5 °C'
$ws.Range("E22").Value = '2026-02-23 06:20:08'
$ws.Range("H22").Value = '''25%'
$ws.Range("E23").Value = '2026-02-23 06:20:10'
$ws.Range("H23").Value = '''26%'
$ws.Range("E24").Value = '2026-02-23 06:20:13'
$ws.Range("N24").Value = '0.2 °C 5:46 TU'
$ws.Range("O24").Value = '2.3 °C'
$ws.Range("E25").Value = '2026-02-23 06:20:16'
$ws.Range("H25").Value = '''29%'
$ws.Range("O25").Value = '3.3 °C'
$ws.Range("E26").Value = '2026-02-23 06:20:18'
$ws.Range("N26").Value = '5.1 °C 5:57 TU'
$ws.Range("E27").Value = '2026-02-23 06:20:21'
$ws.Range("E28").Value = '2026-02-23 06:20:24'
$ws.Range("J28").Value = '1027.3 hPa'
$ws.Range("N28").Value = '1.7 °C 5:57 TU'
$ws.Range("O28").Value = '3.6 °C'
$ws.Range("E29").Value = '2026-02-23 06:20:26'
$ws.Range("K29").Value = '-0.1 MJ/m2'
$ws.Range("N29").Value = '2.9 °C 5:49 TU'
$ws.Range("E30").Value = '2026-02-23 06:20:29'
$ws.Range("O30").Value = '8.1 °C'
$ws.Range("E31").Value = '2026-02-23 06:20:32'
$ws.Range("H31").Value = '''49%'
$ws.Range("J31").Value = '1024.3 hPa'
$ws.Range("O31").Value = '15.1 °C'
$ws.Range("E32").Value = '2026-02-23 06:20:34'
$ws.Range("L32").Value = '5.4 km/h - 283º 5:42 TU'
$ws.Range("O32").Value = '1.3 °C'
$ws.Range("E33").Value = '2026-02-23 06:20:37'
$ws.Range("J33").Value = '1029.7 hPa'
$ws.Range("N33").Value = '0.2 °C 5:36 TU'
$ws.Range("O33").Value = '2.3 °C'
$ws.Range("E34").Value = '2026-02-23 06:20:40'
$ws.Range("H34").Value = '''47%'
$ws.Range("O34").Value = '2.0 °C'
$ws.Range("E35").Value = '2026-02-23 06:20:43'
$ws.Range("J35").Value = '1026.5 hPa'
$ws.Range("O35").Value = '10.1 °C'
$ws.Range("E36").Value = '2026-02-23 06:20:45'
$ws.Range("J36").Value = '1025.3 hPa'
$ws.Range("O36").Value = '6.8 °C'
$ws.Range("E37").Value = '2026-02-23 06:20:48'
$ws.Range("H37").Value = '''80%'
$ws.Range("N37").Value = '1.0 °C 5:57 TU'
$ws.Range("O37").Value = '3.7 °C'
$ws.Range("E38").Value = '2026-02-23 06:20:51'
$ws.Range("H38").Value = '''73%'
$ws.Range("O38").Value = '6.3 °C'
$ws.Range("E39").Value = '2026-02-23 06:20:53'
$ws.Range("H39").Value = '''24%'
$ws.Range("M39").Value = '4.9 °C 5:47 TU'
$ws.Range("E40").Value = '2026-02-23 06:20:56'
$ws.Range("N40").Value = '0.4 °C 5:30 TU'
$ws.Range("O40").Value = '1.9 °C'
$ws.Range("E41").Value = '2026-02-23 06:20:59'
$ws.Range("H41").Value = '''86%'
$ws.Range("N41").Value = '5.8 °C 5:37 TU'
$ws.Range("O41").Value = '7.2 °C'
$ws.Range("E42").Value = '2026-02-23 06:21:02'
$ws.Range("N42").Value = '4.9 °C 5:45 TU'
$ws.Range("O42").Value = '6.0 °C'
$ws.Range("E43").Value = '2026-02-23 06:21:04'
$ws.Range("H43").Value = '''95%'
$ws.Range("K43").Value = '-0.1 MJ/m2'
$ws.Range("N43").Value = '1.7 °C 5:59 TU'
$ws.Range("O43").Value = '3.8 °C'
$ws.Range("E44").Value = '2026-02-23 06:21:07'
$ws.Range("H44").Value = '''37%'
$ws.Range("E45").Value = '2026-02-23 06:21:10'
$ws.Range("H45").Value = '''66%'
$ws.Range("J45").Value = '1030.5 hPa'
$ws.Range("E46").Value = '2026-02-23 06:21:13'
$ws.Range("N46").Value = '0.6 °C 5:59 TU'
$ws.Range("O46").Value = '2.0 °C'
